# Fixed Gantt Chart Sprint 2
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Project Planner")

# Row 5 - "Everyone" / "Play the game"
$ws.Range("E5").Value = 10
$ws.Range("F5").Value = 4
$ws.Range("G5").Value = 10

# Row 6 - "Everyone" / "Make 3 user stories"
$ws.Range("D6").Value = 12
$ws.Range("E6").Value = 4
$ws.Range("F6").Value = 14
$ws.Range("G6").Value = 4

# Update the active selection shown when the sheet is opened
$ws.Activate()
$ws.Range("G7").Select()
